$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.690.13'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '3.814.50'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.59'
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.60'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").Value = '3.812.21'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.69'
$ws.Range("E12").Value = '  +4.83%  '
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.47'
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").Value = '4.454.63'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("D16").Value = '3.812.33'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '68.670.91'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.07'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.11'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '464.15'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.65'
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("E24").Value = '  +2.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.78'
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.02'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.12'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").Value = '3.964.67'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("E31").Value = '  -5.25%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.23'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.00'
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.03'
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  +1.40%  '
$ws.Range("E38").Value = '  +6.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.87'
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.981'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.16'
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '157.07'
$ws.Range("E44").Value = '  +3.50%  '
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.94'
$ws.Range("E46").Value = '  -3.92%  '
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.59'
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.88'
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '378.65'
$ws.Range("E51").Value = '  -2.76%  '
